$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.536.06"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.563.90"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.34"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0885"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "1.786.23"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "1.580.54"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "28.518.78"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.67"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.03%  "
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").Value = "1.398.80"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.787"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.54%  "
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").Value = "1.700.00"
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("E50").Value = "  -4.47%  "
$ws.Range("E51").Value = "  -0.73%  "
